$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the Price/Volume columns for the rows we touch so that
# numeric-looking strings (e.g. "2.17") are stored as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.691.09'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.584.39'
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("E4").Value = '  +1.20%  '
$ws.Range("D5").Value = '207.07'
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").Value = '0.505'
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("E7").Value = '  +1.28%  '
$ws.Range("D8").Value = '22.25'
$ws.Range("E8").Value = '  -3.90%  '
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").Value = '0.0590'
$ws.Range("E10").Value = '  -2.44%  '
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").Value = '1.810.00'
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("D13").Value = '1.599.29'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -2.80%  '
$ws.Range("D16").Value = '27.671.03'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '63.26'
$ws.Range("E17").Value = '  -2.10%  '
$ws.Range("D18").Value = '219.40'
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").Value = '7.34'
$ws.Range("E19").Value = '  -4.08%  '
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("E21").Value = '  +1.13%  '
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").Value = '  -3.77%  '
$ws.Range("E23").Value = '  -5.54%  '
$ws.Range("E24").Value = '  -3.12%  '
$ws.Range("D25").Value = '154.80'
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").Value = '6.86'
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  +1.20%  '
$ws.Range("E28").Value = '  -1.94%  '
$ws.Range("E29").Value = '  -3.12%  '
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("E32").Value = '  -4.03%  '
$ws.Range("D33").Value = '1.382.62'
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("E34").Value = '  -4.82%  '
$ws.Range("E35").Value = '  -4.29%  '
$ws.Range("D36").Value = '0.971'
$ws.Range("E36").Value = '  -3.28%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -2.82%  '
$ws.Range("E39").Value = '  -3.04%  '
$ws.Range("D40").Value = '0.819'
$ws.Range("E40").Value = '  -2.73%  '
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '0.978'
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("D45").Value = '1.74'
$ws.Range("E45").Value = '  -4.01%  '
$ws.Range("E46").Value = '  -2.96%  '
$ws.Range("D47").Value = '1.721.20'
$ws.Range("E47").Value = '  -2.07%  '
$ws.Range("D48").Value = '88.25'
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("E49").Value = '  +6.51%  '
$ws.Range("E50").Value = '  -3.89%  '
$ws.Range("E51").Value = '  -0.42%  '

# Row 43/44: MXToken and Aave swapped rank positions.
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '2.17'
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '63.54'
$ws.Range("E44").Value = '  -2.93%  '

# Restore the original (unformatted) style on the touched cells now that the
# values have been committed as text, so no stray number-format styling remains.
$ws.Range("D2:E51").ClearFormats()

